# "End of year data from Matt" - updated poker stats on Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 - Jon
$ws.Range("D4").Value = 185
$ws.Range("I4").Value = 675
$ws.Range("J4").Value = 3.65

# Row 5 - Maisy
$ws.Range("D5").Value = 113
$ws.Range("I5").Value = 426
$ws.Range("J5").Value = 3.77

# Row 6 - Mark
$ws.Range("D6").Value = 132
$ws.Range("I6").Value = 471
$ws.Range("J6").Value = 3.57

# Row 7 - Matt
$ws.Range("D7").Value = 180
$ws.Range("J7").Value = 3.81

# Row 8 - Pepe
$ws.Range("D8").Value = 95
$ws.Range("I8").Value = 353
$ws.Range("J8").Value = 3.72

# Row 10 - Richard
$ws.Range("D10").Value = 132
$ws.Range("I10").Value = 550

# Row 13 - Andy
$ws.Range("D13").Value = 189
$ws.Range("I13").Value = 781

# Row 14 - Anthony
$ws.Range("D14").Value = 113
$ws.Range("I14").Value = 448
$ws.Range("J14").Value = 3.96
